$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '31.092.97'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.968.12'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +2.84%  '
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.17%  '
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '248.48'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +1.77%  '
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4875'
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +0.74%  '
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '44.76'
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.80%  '
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.2961'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  +2.15%  '
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.06827'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -0.81%  '
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -0.16%  '
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '107.80'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -2.80%  '
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.972.58'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  +2.96%  '
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.07792'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +2.94%  '
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.461'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +1.97%  '
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.7050'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +5.25%  '
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '287.85'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -1.79%  '
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '31.111.22'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +1.83%  '
$c.Style = "Normal"
$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = 'BitDAO'
$c.Style = "Normal"
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.4980'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +13.14%  '
$c.Style = "Normal"
$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = 'Avalanche'
$c.Style = "Normal"
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.24'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +1.75%  '
$c.Style = "Normal"
$c = $ws.Range("B21")
$c.NumberFormat = "@"
$c.Value = 'ShibaInu'
$c.Style = "Normal"
$c = $ws.Range("C21")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.000007748'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  +2.04%  '
$c.Style = "Normal"
$c = $ws.Range("B22")
$c.NumberFormat = "@"
$c.Value = 'WrappedliquidstakedEther2.0'
$c.Style = "Normal"
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.224.01'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  +2.50%  '
$c.Style = "Normal"
$c = $ws.Range("B23")
$c.NumberFormat = "@"
$c.Value = 'Uniswap'
$c.Style = "Normal"
$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.642'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  +2.54%  '
$c.Style = "Normal"
$c = $ws.Range("B24")
$c.NumberFormat = "@"
$c.Value = 'Dai'
$c.Style = "Normal"
$c = $ws.Range("C24")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"
$c = $ws.Range("B25")
$c.NumberFormat = "@"
$c.Value = 'BinanceUSD'
$c.Style = "Normal"
$c = $ws.Range("C25")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +0.29%  '
$c.Style = "Normal"
$c = $ws.Range("B26")
$c.NumberFormat = "@"
$c.Value = 'Chainlink'
$c.Style = "Normal"
$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '6.631'
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +3.66%  '
$c.Style = "Normal"
$c = $ws.Range("B27")
$c.NumberFormat = "@"
$c.Value = 'Cosmos'
$c.Style = "Normal"
$c = $ws.Range("C27")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.03'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +6.08%  '
$c.Style = "Normal"
$c = $ws.Range("B28")
$c.NumberFormat = "@"
$c.Value = 'Monero'
$c.Style = "Normal"
$c = $ws.Range("C28")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '170.47'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +3.50%  '
$c.Style = "Normal"
$c = $ws.Range("B29")
$c.NumberFormat = "@"
$c.Value = 'EthereumClassic'
$c.Style = "Normal"
$c = $ws.Range("C29")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '20.07'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -0.77%  '
$c.Style = "Normal"
$c = $ws.Range("B30")
$c.NumberFormat = "@"
$c.Value = 'LidoDAOToken'
$c.Style = "Normal"
$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.200'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  +5.16%  '
$c.Style = "Normal"
$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = 'Stellar'
$c.Style = "Normal"
$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.1070'
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c.Style = "Normal"
$c = $ws.Range("B32")
$c.NumberFormat = "@"
$c.Value = 'Toncoin'
$c.Style = "Normal"
$c = $ws.Range("C32")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.446'
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  +0.66%  '
$c.Style = "Normal"
$c = $ws.Range("B33")
$c.NumberFormat = "@"
$c.Value = 'Filecoin'
$c.Style = "Normal"
$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.846'
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +19.79%  '
$c.Style = "Normal"
$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = 'InternetComputer(DFINITY)'
$c.Style = "Normal"
$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.528'
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +9.52%  '
$c.Style = "Normal"
$c = $ws.Range("B35")
$c.NumberFormat = "@"
$c.Value = 'Hedera'
$c.Style = "Normal"
$c = $ws.Range("C35")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.05110'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  +2.57%  '
$c.Style = "Normal"
$c = $ws.Range("B36")
$c.NumberFormat = "@"
$c.Value = 'ImmutableX'
$c.Style = "Normal"
$c = $ws.Range("C36")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.7733'
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +5.22%  '
$c.Style = "Normal"
$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = 'ARBITRUM'
$c.Style = "Normal"
$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.172'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +3.51%  '
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.02050'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +1.21%  '
$c.Style = "Normal"
$c = $ws.Range("B39")
$c.NumberFormat = "@"
$c.Value = 'HuobiToken'
$c.Style = "Normal"
$c = $ws.Range("C39")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.735'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +0.93%  '
$c.Style = "Normal"
$c = $ws.Range("B40")
$c.NumberFormat = "@"
$c.Value = 'MXToken'
$c.Style = "Normal"
$c = $ws.Range("C40")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.734'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +1.84%  '
$c.Style = "Normal"
$c = $ws.Range("B41")
$c.NumberFormat = "@"
$c.Value = 'FraxShare'
$c.Style = "Normal"
$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.516'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +11.86%  '
$c.Style = "Normal"
$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = 'RenderToken'
$c.Style = "Normal"
$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.138'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +6.08%  '
$c.Style = "Normal"
$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = 'TrustWalletToken'
$c.Style = "Normal"
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.8907'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +3.31%  '
$c.Style = "Normal"
$c = $ws.Range("B44")
$c.NumberFormat = "@"
$c.Value = 'TheSandbox'
$c.Style = "Normal"
$c = $ws.Range("C44")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.4497'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  +1.41%  '
$c.Style = "Normal"
$c = $ws.Range("B45")
$c.NumberFormat = "@"
$c.Value = 'Quant'
$c.Style = "Normal"
$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '110.02'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +0.44%  '
$c.Style = "Normal"
$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = 'Aave'
$c.Style = "Normal"
$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '72.99'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +5.54%  '
$c.Style = "Normal"
$c = $ws.Range("B47")
$c.NumberFormat = "@"
$c.Value = 'PaxDollar'
$c.Style = "Normal"
$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"
$c = $ws.Range("B48")
$c.NumberFormat = "@"
$c.Value = 'Aptos'
$c.Style = "Normal"
$c = $ws.Range("C48")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.552'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +5.00%  '
$c.Style = "Normal"
$c = $ws.Range("B49")
$c.NumberFormat = "@"
$c.Value = 'Maker'
$c.Style = "Normal"
$c = $ws.Range("C49")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '989.62'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  +17.76%  '
$c.Style = "Normal"
$c = $ws.Range("B50")
$c.NumberFormat = "@"
$c.Value = 'EnergySwap'
$c.Style = "Normal"
$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '9.503'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +2.87%  '
$c.Style = "Normal"
$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = 'Algorand'
$c.Style = "Normal"
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.1270'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +3.62%  '
$c.Style = "Normal"
